# Update parameter names in the Accuracy sheet to the new naming convention.
# Cells are written bottom-to-top (row 19 first, row 2 last) so that the
# workbook's shared-string table regenerates new entries in the same order
# as the target file (E coli, chlorophyll a, orthoP, sp conductivity,
# DO concentration, water temperature).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "E.coli" -> "E coli" (rows 19-22)
$ws.Range("A19").Value = "E coli"
$ws.Range("A20").Value = "E coli"
$ws.Range("A21").Value = "E coli"
$ws.Range("A22").Value = "E coli"

# "Chl-a" -> "chlorophyll a" (rows 17-18)
$ws.Range("A17").Value = "chlorophyll a"
$ws.Range("A18").Value = "chlorophyll a"

# "ortho – P" -> "orthoP" (rows 12-13)
$ws.Range("A12").Value = "orthoP"
$ws.Range("A13").Value = "orthoP"

# "Conductivity" -> "sp conductivity" (rows 6-7)
$ws.Range("A6").Value = "sp conductivity"
$ws.Range("A7").Value = "sp conductivity"

# "DO" -> "DO concentration" (rows 4-5)
$ws.Range("A4").Value = "DO concentration"
$ws.Range("A5").Value = "DO concentration"

# "Water Temp" -> "water temperature" (row 2)
$ws.Range("A2").Value = "water temperature"

# Move the active selection on the frozen (bottom-right) pane to D1,
# matching the saved view state in the target workbook.
$ws.Range("D1").Select()
